# Updated cryptos list (Price / Volume(1h) refresh), matching the
# "Updated cryptos list on Fri Jul 26 16:59:46 UTC 2024 with GitHub Actions"
# commit. Numeric-looking text values (e.g. "1.00", "31.54") are written
# with the cell pre-formatted as Text ("@") so Excel keeps them as strings
# instead of silently converting them to numbers; the style is then reset
# to "Normal" so no stray number-format/quote-prefix is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.461.32'
$ws.Range('E2').Value = '  +3.97%  '
$ws.Range('D3').Value = '3.248.95'
$ws.Range('E3').Value = '  +2.90%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.07%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -3.00%  '
$ws.Range('D9').Value = '3.251.57'
$ws.Range('E9').Value = '  +2.99%  '
$ws.Range('E10').Value = '  +5.09%  '
$ws.Range('E11').Value = '  +3.09%  '
$ws.Range('E12').Value = '  +5.18%  '
$ws.Range('D13').Value = '3.814.63'
$ws.Range('E13').Value = '  +3.00%  '
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.67'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.75%  '
$ws.Range('D16').Value = '67.474.29'
$ws.Range('E16').Value = '  +4.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000167'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.92%  '
$ws.Range('D18').Value = '3.254.70'
$ws.Range('E18').Value = '  +3.02%  '
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '376.55'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.63'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('E24').Value = '  +3.65%  '
$ws.Range('E26').Value = '  +1.85%  '
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('E28').Value = '  +2.58%  '
$ws.Range('E29').Value = '  +0.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.79'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.83%  '
$ws.Range('E31').Value = '  +3.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '22.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.15%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.27'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.68%  '
$ws.Range('E35').Value = '  +4.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '164.03'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.40%  '
$ws.Range('E37').Value = '  +3.85%  '
$ws.Range('E38').Value = '  +1.99%  '
$ws.Range('E39').Value = '  +5.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.78'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +12.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.58'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.61'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '358.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.15%  '
$ws.Range('D45').Value = '2.726.86'
$ws.Range('E45').Value = '  +2.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.32'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.87'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.85%  '
$ws.Range('E48').Value = '  +3.34%  '
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('E50').Value = '  +6.83%  '
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.54'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.36%  '
